# Add "RCI" as a new header in column F (row 1) and populate RCI values
# for each scale/index row (2-13), mirroring the pre/post T-score and
# percentile-rank columns already present in the brief.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "RCI"

$rciValues = @(
    0,
    0.97,
    -0.72,
    -0.5,
    -0.73,
    -0.32,
    0,
    -0.45,
    0,
    -1.09,
    -0.72,
    -1.22
)

for ($i = 0; $i -lt $rciValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $rciValues[$i]
}

# Give the new column a sensible width, matching the sizing used for the
# other data columns on this sheet.
$ws.Columns.Item(6).ColumnWidth = 17.5

# Move / update the active selection like the author's saved state.
$ws.Range("F14").Select()
